# "Generate Report for Handoff"
# e2e\b.md has a fresh handoff xliff generated; mark it "Ready for handoff"
# on the Overview + per-locale sheets, and flag the stale handback files.

$wb = $excel.ActiveWorkbook

$status        = "Ready for handoff"
$genDate       = "2016-10-17 14:25:52"
$zhHandoff     = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$deHandoff     = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$zhHandoffDate = "2016-10-17 14:25:30"
$deHandoffDate = "2016-10-17 14:25:52"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1cd2f7a7df83cdffb78a3e48d41ad44709e80b8b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5bf40b5f1c17bc28134d05e29f2f2ac8710ad21/e2e/b.md."

# --- Overview sheet: row for b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $genDate

# --- zh-cn sheet: row for b.md (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $status
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G3").Value = $zhHandoff
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666

# --- de-de sheet: row for b.md (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $status
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("F3").Style = "Normal"
$wsDe.Range("G3").Value = $deHandoff
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666
